$wb = $excel.ActiveWorkbook

# ---- Metadata sheet ----
$meta = $wb.Worksheets.Item("Metadata")

# Version: 2.2.0-ballot -> 2.1.0
$meta.Range("B3").Value = "2.1.0"

# Date: 2025-12-19T08:32:44+00:00 -> 2025-12-19T08:44:55+00:00
$meta.Range("B8").Value = "2025-12-19T08:44:55+00:00"

# Base Definition: drop the trailing "|2.1.0" version pin
$meta.Range("B18").Value = "https://hl7.fr/ig/fhir/core/StructureDefinition/fr-core-human-name"

# ---- Elements sheet ----
$el = $wb.Worksheets.Item("Elements")

# Type(s) for assemblyOrder extension: drop the trailing "|5.2.0" version pin
$el.Range("K5").Value = "Extension {humanname-assembly-order}`n"

# Binding Value Set for HumanName.use: drop the trailing "|4.0.1" version pin
$el.Range("Z6").Value = "http://hl7.org/fhir/ValueSet/name-use"

# Binding Value Set for HumanName.prefix: drop the trailing "|20230331120000" version pin
$el.Range("Z10").Value = "https://mos.esante.gouv.fr/NOS/JDV_J245-Civilite-CISIS/FHIR/JDV-J245-Civilite-CISIS"

# Binding Value Set for HumanName.suffix: drop the trailing "|20200424120000" version pin
$el.Range("Z11").Value = "https://mos.esante.gouv.fr/NOS/JDV_J79-CiviliteExercice-RASS/FHIR/JDV-J79-CiviliteExercice-RASS"

# The shortened text above narrows the best-fit autosized width of columns K (Type(s))
# and Z (Binding Value Set); reflect the resulting column widths.
$el.Columns.Item(11).ColumnWidth = 33.0
$el.Columns.Item(26).ColumnWidth = 76.5
